$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E29) is reordered to list periods most-recent-first
# (2105 down to 2004) instead of oldest-first (2004 up to 2105).
$periodos = @("2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004")

$row = 16
foreach ($periodo in $periodos) {
    $ws.Range("E$row").Value = $periodo
    $row = $row + 1
}

# The "Valor Mora" figures tied to the first/last period rows swap along with
# the reordering: row 16 (now period 2105) takes the value that used to sit
# on row 29, and row 29 (now period 2004) takes the value that used to sit
# on row 16.
$ws.Range("F16").Value = 28090
$ws.Range("F29").Value = 26919
